# Auto-generated edit script: update cryptos list values for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.246.68'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '1.915.11'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '''0.7401'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").Value = '''243.92'
$ws.Range("E6").Value = '  -2.36%  '
$ws.Range("D7").Value = '''1.002'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '''0.3131'
$ws.Range("E8").Value = '  -2.63%  '
$ws.Range("D9").Value = '''27.00'
$ws.Range("E9").Value = '  -3.64%  '
$ws.Range("D10").Value = '''0.06949'
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '''0.7727'
$ws.Range("E11").Value = '  -2.17%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.974.23'
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.07972'
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").Value = '''5.263'
$ws.Range("E14").Value = '  -2.29%  '
$ws.Range("D15").Value = '''91.32'
$ws.Range("E15").Value = '  -3.35%  '
$ws.Range("D16").Value = '30.288.01'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '''14.16'
$ws.Range("E17").Value = '  -2.97%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").Value = '''245.51'
$ws.Range("E18").Value = '  -2.85%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '''5.776'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").Value = '''0.000007823'
$ws.Range("E20").Value = '  -2.72%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.187.45'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '''1.003'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '''1.002'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = '''6.637'
$ws.Range("E24").Value = '  -2.80%  '
$ws.Range("D25").Value = '''9.386'
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = '''165.24'
$ws.Range("E26").Value = '  +0.45%  '
$ws.Range("D27").Value = '''18.95'
$ws.Range("E27").Value = '  -0.86%  '
$ws.Range("D28").Value = '''0.1269'
$ws.Range("E28").Value = '  -5.37%  '
$ws.Range("D29").Value = '''2.130'
$ws.Range("E29").Value = '  -8.20%  '
$ws.Range("D30").Value = '''1.365'
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("D31").Value = '''1.546'
$ws.Range("E31").Value = '  +0.65%  '
$ws.Range("D32").Value = '''4.313'
$ws.Range("E32").Value = '  -2.70%  '
$ws.Range("D33").Value = '''4.071'
$ws.Range("E33").Value = '  -1.96%  '
$ws.Range("D34").Value = '''0.05162'
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("D35").Value = '''1.293'
$ws.Range("E35").Value = '  +0.25%  '
$ws.Range("D36").Value = '''0.7464'
$ws.Range("E36").Value = '  -0.50%  '
$ws.Range("D37").Value = '''2.773'
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").Value = '''0.01933'
$ws.Range("E38").Value = '  -1.81%  '
$ws.Range("D39").Value = '''2.785'
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").Value = '''6.353'
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("D41").Value = '''76.00'
$ws.Range("E41").Value = '  -2.64%  '
$ws.Range("D42").Value = '''0.4466'
$ws.Range("E42").Value = '  -1.27%  '
$ws.Range("D43").Value = '''1.944'
$ws.Range("E43").Value = '  -2.66%  '
$ws.Range("D44").Value = '''1.002'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").Value = '''0.8348'
$ws.Range("E45").Value = '  -0.84%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '''7.667'
$ws.Range("E46").Value = '  +1.47%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '''101.44'
$ws.Range("E47").Value = '  -0.86%  '
$ws.Range("D48").Value = '''9.854'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''36.91'
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.1213'
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '''942.25'
$ws.Range("E51").Value = '  -4.54%  '

# Reset style on cells that were forced to text via leading apostrophe,
# so no stray quotePrefix formatting remains on the cell (matches original formatting).
$forceCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D17","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $forceCells) {
    $ws.Range($addr).Style = "Normal"
}
